# Commit: "Change names from *img to img*"
# Rename the seven *img sheets to img* (prefix moved from the end to the
# front of the name), and make the newly-reordered "imge" sheet (formerly
# "eimg", now the 17th / last tab) the active sheet instead of "ebday".

$wb = $excel.ActiveWorkbook

$renames = @{
    "himg" = "imgh";
    "timg" = "imgt";
    "simg" = "imgs";
    "gimg" = "imgg";
    "wimg" = "imgw";
    "bimg" = "imgb";
    "eimg" = "imge";
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# Move the active tab from "ebday" (index 5) to "imge" (the renamed
# "eimg", now at index 16 / the last sheet).
$wb.Worksheets.Item("imge").Activate()
